$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header row (now "Input" / "Output" / "Name" digraph-style headers)
$ws.Range("A1").Value = "Input"
$ws.Range("B1").Value = "Output"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Data (bytes)"
$ws.Range("E1").Value = "True Pass"
$ws.Range("F1").Value = "True Discard"
$ws.Range("G1").Value = "False Pass (alpha)"
$ws.Range("H1").Value = "False Discard (beta)"
$ws.Range("I1").Value = "Compression"
$ws.Range("J1").Value = "PU 140"
$ws.Range("K1").Value = "PU 200"

# Row 2 - Tracking
$ws.Range("A2").Value = "Tracking"
$ws.Range("B2").Value = "Intermediate"
$ws.Range("C2").Value = "Tracking"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").ClearContents()

# Row 3 - Timing
$ws.Range("A3").Value = "Timing"
$ws.Range("B3").Value = "Intermediate"
$ws.Range("C3").Value = "Timing"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

# Row 4 - Calorimetry (new row)
$ws.Range("A4").Value = "Calorimetry"
$ws.Range("B4").Value = "Intermediate"
$ws.Range("C4").Value = "Calorimetry"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0

# Row 5 - Muon (new row)
$ws.Range("A5").Value = "Muon"
$ws.Range("B5").Value = "Intermediate"
$ws.Range("C5").Value = "Muon"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0

# Row 6 - Intermediate -> Global : Level-1 Trigger (previously row 2)
$ws.Range("A6").Value = "Intermediate"
$ws.Range("B6").Value = "Global"
$ws.Range("C6").Value = "Level-1 Trigger"
$ws.Range("D6").Formula = "=INT(1000000*K6)"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 400
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0.26
$ws.Range("K6").Value = 0.26

# Row 7 - Global -> Disk : High-Level Trigger (previously row 3)
$ws.Range("A7").Value = "Global"
$ws.Range("B7").Value = "Disk"
$ws.Range("C7").Value = "High-Level Trigger"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 20
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

[void]$ws.Range("C16").Select()
